# X-Helg 2017 Statistikk - update for 5. Des (funn på publiseringsdato)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- 1. Extend formatting of new rows 20:22 to match existing data rows (copy format from row 19) ---
$ws.Range("A19:G19").Copy() | Out-Null
$ws.Range("A20:G22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Helper to write a text cell (nick/plass/ftf/etc columns, all Text-formatted) ---
function Set-Text($addr, $val) {
    if ($null -eq $val -or $val -eq "") {
        $ws.Range($addr).ClearContents() | Out-Null
    } else {
        $ws.Range($addr).Value = $val
    }
}

# --- 3. Helper to write the numeric Total column (G), which is formatted as Text (@) so the
#        normal .Value setter would coerce the number into a brand new shared string instead of
#        a numeric cell. Flip the format to General while writing, then restore Text afterwards
#        so the stored value is a real number but still looks/behaves like the original cells. ---
function Set-Total($addr, $val) {
    $fmt = $ws.Range($addr).NumberFormat
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = $fmt
}

# --- 4. Row data: Plass(A), Nick(B), Publisert selv(C), FTF(D), Funn på publ dato(E), Funn i Desember(F), Total(G)
$rows = @(
    @{ Row=4;  A='1'; B='siljejandersen';  C=$null;            D='[#3]* [#2]*  (4)'; E='[#1] [#2] [#3]  (6)'; F=$null;                 G=10 },
    @{ Row=5;  A='2'; B='johs1988';        C=$null;            D='[#3]* [#2]*  (4)'; E='[#1]  (2)';            F='[#2] [#3]  (2)';     G=8  },
    @{ Row=6;  A='3'; B='O-K Haukland';    C=$null;            D='[#2]*  (2)';       E='[#2]  (2)';            F='[#1] [#3]  (2)';     G=6  },
    @{ Row=7;  A='4'; B='minni09';         C='[#2]  (3)';      D=$null;              E=$null;                  F='[#1] [#3]  (2)';     G=5  },
    @{ Row=8;  A='5'; B='TeamCOR';         C=$null;            D=$null;              E='[#1] [#2]  (4)';      F=$null;                 G=4  },
    @{ Row=9;  A='5'; B='TeamLadybug<3<3'; C=$null;            D='[#1]  (3)';        E=$null;                  F='[#1]  (1)';          G=4  },
    @{ Row=10; A='5'; B='bleikfis';        C='[#1]  (3)';      D=$null;              E=$null;                  F='[#2]  (1)';          G=4  },
    @{ Row=11; A='6'; B='cara2006';        C=$null;            D=$null;              E='[#1]  (2)';            F='[#2]  (1)';          G=3  },
    @{ Row=12; A='6'; B='SisselHultgreen'; C='[#13]  (3)';     D=$null;              E=$null;                  F=$null;                 G=3  },
    @{ Row=13; A='6'; B='GunnarKolskog';   C='[#3]  (3)';      D=$null;              E=$null;                  F=$null;                 G=3  },
    @{ Row=14; A='6'; B='Team Lynis';      C='[#4]  (3)';      D=$null;              E=$null;                  F=$null;                 G=3  },
    @{ Row=15; A='6'; B='TeamPolhøgda';    C='[#5]  (3)';      D=$null;              E=$null;                  F=$null;                 G=3  },
    @{ Row=16; A='7'; B='Kransa';          C=$null;            D=$null;              E=$null;                  F='[#1] [#2]  (2)';     G=2  },
    @{ Row=17; A='7'; B='Onyx Black';      C=$null;            D=$null;              E='[#2]  (2)';            F=$null;                 G=2  },
    @{ Row=18; A='7'; B='SonjaJ';          C=$null;            D=$null;              E='[#2]  (2)';            F=$null;                 G=2  },
    @{ Row=19; A='7'; B='silyam';          C=$null;            D=$null;              E=$null;                  F='[#1] [#2]  (2)';     G=2  },
    @{ Row=20; A='7'; B='tomnor';          C=$null;            D=$null;              E='[#2]  (2)';            F=$null;                 G=2  },
    @{ Row=21; A='8'; B='annesto';         C=$null;            D=$null;              E=$null;                  F='[#2]  (1)';          G=1  },
    @{ Row=22; A='8'; B='dogteam';         C=$null;            D=$null;              E=$null;                  F='[#2]  (1)';          G=1  }
)

foreach ($r in $rows) {
    $row = $r.Row
    Set-Text "A$row" $r.A
    Set-Text "B$row" $r.B
    Set-Text "C$row" $r.C
    Set-Text "D$row" $r.D
    Set-Text "E$row" $r.E
    Set-Text "F$row" $r.F
    Set-Total "G$row" $r.G
}

# --- 5. Header banner text (merged F2:G2) ---
Set-Text "F2" "26 besøk hittil! 19 team deltok! "

# --- 6. Defined names: extend the statistikk range and the print area to cover the new rows ---
$wb.Names.Item("Ark1!statistikk").RefersTo = "='Ark1'!`$A`$3:`$G`$22"
$wb.Names.Item("Ark1!Print_Area").RefersTo = "='Ark1'!`$A`$1:`$G`$25"
